$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (row 1) text relabeling ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Data value updates (column C: GDP; column AL: Colony flag) ---
$ws.Range("C2").Value = 2870.311589353206
$ws.Range("C3").Value = 5191.140356354663
$ws.Range("AL3").Value = 1
$ws.Range("C4").Value = 8947.741473873051
$ws.Range("C5").Value = 10594.98659239237
$ws.Range("C6").Value = 1909.084588129339
$ws.Range("C7").Value = 9502.243585046588
$ws.Range("C8").Value = 12227.21453003286
$ws.Range("C9").Value = 6128.19547247793
$ws.Range("C10").Value = 4729.735976516416
$ws.Range("C11").Value = 11155.84524560499
$ws.Range("C12").Value = 14239.03920301361
$ws.Range("C13").Value = 2100.656463590606
$ws.Range("C14").Value = 1286.515571617672
$ws.Range("C15").Value = 2812.435974421079
$ws.Range("C16").Value = 17288.8595992193
$ws.Range("C17").Value = 2898.942214704482
$ws.Range("C18").Value = 5555.389721901988
$ws.Range("AL18").Value = 1
$ws.Range("C19").Value = 9271.398233246389
$ws.Range("C20").Value = 13825.35808833117
$ws.Range("C21").Value = 5082.354756663512
$ws.Range("C22").Value = 6336.709213679884
$ws.Range("C23").Value = 12808.034586422
$ws.Range("C24").Value = 11286.24301624575
$ws.Range("C25").Value = 1955.461557360978
$ws.Range("C26").Value = 10385.96443195552
$ws.Range("C27").Value = 2217.474008566157
$ws.Range("C28").Value = 1357.563719132622
$ws.Range("C29").Value = 1303.425880277445
$ws.Range("C30").Value = 2828.483778716848
$ws.Range("C31").Value = 17610.30663334184
$ws.Range("C32").Value = 2965.153206179127
$ws.Range("C33").Value = 9477.887185090232
$ws.Range("C34").Value = 5660.517066940175
$ws.Range("AL34").Value = 1
$ws.Range("C35").Value = 10883.31535948899
$ws.Range("C36").Value = 2024.117324382548
$ws.Range("C37").Value = 11627.81065059172
$ws.Range("C38").Value = 6711.616186806423
$ws.Range("C39").Value = 5360.226632400601
$ws.Range("C40").Value = 2264.394087033834
$ws.Range("C41").Value = 1410.426304742003
$ws.Range("C42").Value = 2286.013198234259
$ws.Range("C43").Value = 1401.753174264641
$ws.Range("C44").Value = 2612.856880840196
$ws.Range("C45").Value = 3137.260298393558
$ws.Range("C46").Value = 16146.07242861928
$ws.Range("C47").Value = 1640.18070024053
$ws.Range("C48").Value = 9839.050190896
$ws.Range("C49").Value = 5996.49696468919
$ws.Range("C51").Value = 6301.696269820412
$ws.Range("AL51").Value = 1
$ws.Range("C52").Value = 2361.056581219794
$ws.Range("C53").Value = 1441.783971398429
$ws.Range("C54").Value = 2735.187532014817
$ws.Range("C55").Value = 3210.869677115934
$ws.Range("C56").Value = 14093.81249338665
$ws.Range("C57").Value = 1751.664428859304
$ws.Range("C58").Value = 10037.20149040966
$ws.Range("C59").Value = 6114.227214287786
$ws.Range("C61").Value = 6661.86504232374
$ws.Range("AL61").Value = 1
$ws.Range("C62").Value = 10239.48134799327
$ws.Range("C63").Value = 2425.561644739583
$ws.Range("C64").Value = 2886.897484630703
$ws.Range("C65").Value = 3242.636921959078
$ws.Range("C66").Value = 12358.30403621203
$ws.Range("C67").Value = 1875.732161108182
$ws.Range("C68").Value = 10205.79575322194
$ws.Range("C69").Value = 6262.368904654469
$ws.Range("C71").Value = 7026.178156858586
$ws.Range("AL71").Value = 1
